$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Forecast Comparison" - shift the Week_Start_Date column forward by
# one week and update the MyForecast values (Penalty/Reward system rerun).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$weekDates = @{
    2  = "2025-01-12"
    3  = "2025-01-19"
    4  = "2025-01-26"
    5  = "2025-02-02"
    6  = "2025-02-09"
    7  = "2025-02-16"
    8  = "2025-02-23"
    9  = "2025-03-02"
    10 = "2025-03-09"
    11 = "2025-03-16"
    12 = "2025-03-23"
    13 = "2025-03-30"
    14 = "2025-04-06"
    15 = "2025-04-13"
    16 = "2025-04-20"
    17 = "2025-04-27"
}

$myForecast = @{
    2  = 89
    3  = 98
    4  = 101
    5  = 100
    6  = 98
    7  = 99
    8  = 101
    9  = 102
    10 = 104
    11 = 106
    12 = 107
    13 = 105
    14 = 102
    15 = 103
    16 = 110
    17 = 74
}

foreach ($row in 2..17) {
    $cell = $ws1.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $weekDates[$row]

    $ws1.Cells.Item($row, 4).Value = $myForecast[$row]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Summary" - refresh the derived statistics to match the new
# forecast numbers above.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Cells.Item(2, 2).NumberFormat = "@"
$ws2.Cells.Item(2, 2).Value = "2023-02-12 to 2025-01-05"

$ws2.Cells.Item(4, 2).NumberFormat = "@"
$ws2.Cells.Item(4, 2).Value = "171"

$ws2.Cells.Item(5, 2).NumberFormat = "@"
$ws2.Cells.Item(5, 2).Value = "60"

$ws2.Cells.Item(7, 2).NumberFormat = "@"
$ws2.Cells.Item(7, 2).Value = "36"

$ws2.Cells.Item(8, 2).NumberFormat = "@"
$ws2.Cells.Item(8, 2).Value = "4226 units"

$ws2.Cells.Item(9, 2).NumberFormat = "@"
$ws2.Cells.Item(9, 2).Value = "1599"

$ws2.Cells.Item(10, 2).NumberFormat = "@"
$ws2.Cells.Item(10, 2).Value = "788"

$ws2.Cells.Item(11, 2).NumberFormat = "@"
$ws2.Cells.Item(11, 2).Value = "388"

$ws2.Cells.Item(12, 2).NumberFormat = "@"
$ws2.Cells.Item(12, 2).Value = "110"

$ws2.Cells.Item(13, 2).NumberFormat = "@"
$ws2.Cells.Item(13, 2).Value = "2025-04-20"

$ws2.Cells.Item(14, 2).NumberFormat = "@"
$ws2.Cells.Item(14, 2).Value = "74"

$ws2.Cells.Item(15, 2).NumberFormat = "@"
$ws2.Cells.Item(15, 2).Value = "2025-04-27"
